$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the "Price" (D) column hold text values formatted like numbers
# (e.g. thousand-separated "70.724.48"). Force text format so COM does not
# coerce them into actual floating point numbers and lose their exact text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.724.48"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "3.527.39"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "624.28"
$ws.Range("E5").Value = "  +2.67%  "

$ws.Range("D6").Value = "173.62"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("D8").Value = "3.522.54"
$ws.Range("E8").Value = "  -2.01%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "0.198"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").Value = "7.01"
$ws.Range("E11").Value = "  -6.50%  "

$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  -1.84%  "

$ws.Range("D13").Value = "46.43"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").Value = "0.0000277"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "4.096.33"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").Value = "8.40"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "609.75"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "3.528.26"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").Value = "70.794.86"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").Value = "17.74"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("D22").Value = "0.884"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").Value = "9.13"
$ws.Range("E23").Value = "  -3.13%  "

$ws.Range("D24").Value = "98.57"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("E25").Value = "  -4.29%  "

$ws.Range("D26").Value = "3.76"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "2.57"
$ws.Range("E28").Value = "  -3.88%  "

$ws.Range("D29").Value = "33.75"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  -3.47%  "

$ws.Range("D31").Value = "3.01"
$ws.Range("E31").Value = "  -2.92%  "

$ws.Range("E32").Value = "  -5.56%  "

$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").Value = "631.31"
$ws.Range("E34").Value = "  -2.04%  "

$ws.Range("D35").Value = "6.79"
$ws.Range("E35").Value = "  -5.35%  "

$ws.Range("E36").Value = "  -3.17%  "

$ws.Range("D37").Value = "10.79"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").Value = "0.0474"
$ws.Range("E38").Value = "  -2.85%  "

$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  -10.49%  "

$ws.Range("D40").Value = "56.96"
$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("D43").Value = "3.349.98"
$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("D44").Value = "0.0₃0722"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").Value = "2.96"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("E46").Value = "  -4.74%  "

$ws.Range("D47").Value = "31.89"
$ws.Range("E47").Value = "  -4.26%  "

$ws.Range("D48").Value = "2.54"
$ws.Range("E48").Value = "  -5.64%  "

$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("D50").Value = "133.76"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.01%  "

